$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns appended after the existing "Landmark" column (L):
# BU, UT, Feeder, DT, Tariff, "Metered  Status" -> M1:R1
$ws.Range("M1:R1").Font.Bold = $true
$ws.Range("M1").Value = "BU"
$ws.Range("N1").Value = "UT"
$ws.Range("O1").Value = "Feeder"
$ws.Range("P1").Value = "DT"
$ws.Range("Q1").Value = "Tariff"
$ws.Range("R1").Value = "Metered  Status"

# Column R (18) gets a custom width, like the other header columns.
$ws.Columns.Item(18).ColumnWidth = 18.83

# Reset the lingering "D4" selection left over from editing, back to A1.
[void]$ws.Range("A1").Select()
